$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 505638.28
$ws.Range("I28").Value = 1010328.06
$ws.Range("J28").Value = 948.4545000000001
$ws.Range("K28").Value = 1010328.06
$ws.Range("L28").Value = 948.4545000000001
$ws.Range("M28").Value = -1009843.06
$ws.Range("N28").Value = -1918.4545
$ws.Range("H32").Value = 1338.5714
$ws.Range("J32").Value = 1174
$ws.Range("L32").Value = 1174
$ws.Range("N32").Value = -1826
$ws.Range("H51").Value = 2785.3076
$ws.Range("I51").Value = 2521.8
$ws.Range("J51").Value = 2950
$ws.Range("K51").Value = 2521.8
$ws.Range("L51").Value = 2950
$ws.Range("M51").Value = -2037.8
$ws.Range("N51").Value = -3918
$ws.Range("H58").Value = 1362.125
$ws.Range("J58").Value = 1802.8334
$ws.Range("L58").Value = 5408.5002
$ws.Range("N58").Value = -5708.5002
$ws.Range("H129").Value = 990.1053000000001
$ws.Range("I129").Value = 451.5
$ws.Range("J129").Value = 1053.4706
$ws.Range("K129").Value = 1354.5
$ws.Range("L129").Value = 3160.4118
$ws.Range("M129").Value = 3645.5
$ws.Range("N129").Value = -13160.4118
$ws.Range("H132").Value = 219912.55
$ws.Range("I132").Value = 253679.36
$ws.Range("J132").Value = 39822.89
$ws.Range("K132").Value = 761038.08
$ws.Range("L132").Value = 119468.67
$ws.Range("M132").Value = -758508.08
$ws.Range("N132").Value = -124528.67
$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -50120
$ws.Range("H135").Value = 1382.1111
$ws.Range("I135").Value = 1462.7142
$ws.Range("J135").Value = 1100
$ws.Range("K135").Value = 13164.4278
$ws.Range("L135").Value = 9900
$ws.Range("M135").Value = -10629.4278
$ws.Range("N135").Value = -14970
$ws.Range("H136").Value = 43499.668
$ws.Range("J136").Value = 43499.668
$ws.Range("L136").Value = 43499.668
$ws.Range("N136").Value = -53699.668
$ws.Range("H138").Value = 6495826.5
$ws.Range("I138").Value = 2597.4348
$ws.Range("K138").Value = 7792.3044
$ws.Range("M138").Value = -2652.3044

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1803.0435
$ws.Range("I45").Value = 1308.0952
$ws.Range("J45").Value = 7000
$ws.Range("K45").Value = 1308.0952
$ws.Range("L45").Value = 7000
$ws.Range("M45").Value = -931.0952
$ws.Range("N45").Value = -7754
$ws.Range("H88").Value = 20000
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 20000
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H97").Value = 13333819
$ws.Range("I97").Value = 22222624
$ws.Range("J97").Value = 612.2
$ws.Range("K97").Value = 22222624
$ws.Range("L97").Value = 612.2
$ws.Range("M97").Value = -22222128
$ws.Range("N97").Value = -1604.2
$ws.Range("H122").Value = 8752.714
$ws.Range("I122").Value = 10683.818
$ws.Range("J122").Value = 1672
$ws.Range("K122").Value = 32051.454
$ws.Range("L122").Value = 5016
$ws.Range("M122").Value = -29601.454
$ws.Range("N122").Value = -9916
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 700.5769
$ws.Range("I80").Value = 416.5
$ws.Range("J80").Value = 826.8333
$ws.Range("K80").Value = 416.5
$ws.Range("L80").Value = 826.8333
$ws.Range("M80").Value = 581.5
$ws.Range("N80").Value = -2822.8333
$ws.Range("H83").Value = 700.5769
$ws.Range("I83").Value = 416.5
$ws.Range("J83").Value = 826.8333
$ws.Range("K83").Value = 2082.5
$ws.Range("L83").Value = 4134.1665
$ws.Range("M83").Value = 2909.5
$ws.Range("N83").Value = -14118.1665
$ws.Range("H105").Value = 13336103
$ws.Range("I105").Value = 16669353
$ws.Range("J105").Value = 3100
$ws.Range("K105").Value = 16669353
$ws.Range("L105").Value = 3100
$ws.Range("M105").Value = -16667606
$ws.Range("N105").Value = -6594

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1319.6
$ws.Range("I31").Value = 925.44446
$ws.Range("J31").Value = 1826.3715
$ws.Range("K31").Value = 925.44446
$ws.Range("L31").Value = 1826.3715
$ws.Range("M31").Value = -630.44446
$ws.Range("N31").Value = -2416.3715
$ws.Range("H34").Value = 1319.6
$ws.Range("I34").Value = 925.44446
$ws.Range("J34").Value = 1826.3715
$ws.Range("K34").Value = 925.44446
$ws.Range("L34").Value = 1826.3715
$ws.Range("M34").Value = -723.44446
$ws.Range("N34").Value = -2230.3715
$ws.Range("H132").Value = 1681.8136
$ws.Range("I132").Value = 1126
$ws.Range("K132").Value = 3378
$ws.Range("M132").Value = -848

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 8225.621999999999
$ws.Range("J39").Value = 8225.621999999999
$ws.Range("L39").Value = 24676.866
$ws.Range("N39").Value = -25264.866
$ws.Range("H86").Value = 292.30768
$ws.Range("I86").Value = 220.1
$ws.Range("J86").Value = 533
$ws.Range("K86").Value = 660.3
$ws.Range("L86").Value = 1599
$ws.Range("M86").Value = 525.7
$ws.Range("N86").Value = -3971
$ws.Range("H89").Value = 292.30768
$ws.Range("I89").Value = 220.1
$ws.Range("J89").Value = 533
$ws.Range("K89").Value = 1980.9
$ws.Range("L89").Value = 4797
$ws.Range("M89").Value = 3947.1
$ws.Range("N89").Value = -16653
$ws.Range("H129").Value = 1376.6666
$ws.Range("I129").Value = 538
$ws.Range("J129").Value = 2425
$ws.Range("K129").Value = 1614
$ws.Range("L129").Value = 7275
$ws.Range("M129").Value = 3386
$ws.Range("N129").Value = -17275
$ws.Range("H131").Value = 2395.862
$ws.Range("J131").Value = 2606.1794
$ws.Range("L131").Value = 7818.5382
$ws.Range("N131").Value = -17898.5382
$ws.Range("H132").Value = 963.1786
$ws.Range("I132").Value = 741.7778
$ws.Range("J132").Value = 1068.0526
$ws.Range("K132").Value = 6676.000199999999
$ws.Range("L132").Value = 9612.473399999999
$ws.Range("M132").Value = -4146.000199999999
$ws.Range("N132").Value = -14672.4734

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 727.5143
$ws.Range("I97").Value = 611.5263
$ws.Range("J97").Value = 865.25
$ws.Range("K97").Value = 611.5263
$ws.Range("L97").Value = 865.25
$ws.Range("M97").Value = -115.5263
$ws.Range("N97").Value = -1857.25
$ws.Range("H122").Value = 1011355.4
$ws.Range("I122").Value = 5556005.5
$ws.Range("J122").Value = 1433.1111
$ws.Range("K122").Value = 16668016.5
$ws.Range("L122").Value = 4299.3333
$ws.Range("M122").Value = -16665566.5
$ws.Range("N122").Value = -9199.3333

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 882
$ws.Range("I22").Value = 931.25
$ws.Range("J22").Value = 783.5
$ws.Range("K22").Value = 931.25
$ws.Range("L22").Value = 783.5
$ws.Range("M22").Value = -636.25
$ws.Range("N22").Value = -1373.5
$ws.Range("H27").Value = 882
$ws.Range("I27").Value = 931.25
$ws.Range("J27").Value = 783.5
$ws.Range("K27").Value = 931.25
$ws.Range("L27").Value = 783.5
$ws.Range("M27").Value = -824.25
$ws.Range("N27").Value = -997.5
$ws.Range("H61").Value = 6755.875
$ws.Range("I61").Value = 7221.143
$ws.Range("J61").Value = 3499
$ws.Range("K61").Value = 7221.143
$ws.Range("L61").Value = 3499
$ws.Range("M61").Value = -7019.143
$ws.Range("N61").Value = -3903
$ws.Range("H113").Value = 6755.875
$ws.Range("I113").Value = 7221.143
$ws.Range("J113").Value = 3499
$ws.Range("K113").Value = 7221.143
$ws.Range("L113").Value = 3499
$ws.Range("M113").Value = -5051.143
$ws.Range("N113").Value = -7839
$ws.Range("H122").Value = 3558.8333
$ws.Range("I122").Value = 2676
$ws.Range("J122").Value = 3811.0715
$ws.Range("K122").Value = 8028
$ws.Range("L122").Value = 11433.2145
$ws.Range("M122").Value = -5578
$ws.Range("N122").Value = -16333.2145

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 101440.4
$ws.Range("I122").Value = 144143.42
$ws.Range("K122").Value = 432430.26
$ws.Range("M122").Value = -429980.26
$ws.Range("H125").Value = 33990.91
$ws.Range("J125").Value = 33990.91
$ws.Range("L125").Value = 33990.91
$ws.Range("N125").Value = -43830.91
$ws.Range("H126").Value = 101090.5
$ws.Range("I126").Value = 112100.555
$ws.Range("K126").Value = 336301.665
$ws.Range("M126").Value = -333831.665
